# Add a new "Nam" (Year) summary sheet after the four quarter sheets,
# pulling together the headers into one combined view, and make it active.

$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add()
$ws.Move($null, $lastSheet)
$ws.Name = "Nam"

# Header row 1: quarter labels, merged across KH/TH pairs.
$ws.Range("A1").Value = "Chỉ tiêu"
$ws.Range("B1").Value = "Quý I"
$ws.Range("D1").Value = "Quý II"
$ws.Range("F1").Value = "Quý III"
$ws.Range("H1").Value = "Quý IV"
$ws.Range("J1").Value = "Tỏng"

# Header row 2: KH/TH sub-labels under each quarter.
$ws.Range("B2").Value = "KH"
$ws.Range("C2").Value = "TH"
$ws.Range("D2").Value = "KH"
$ws.Range("E2").Value = "TH"
$ws.Range("F2").Value = "KH"
$ws.Range("G2").Value = "TH"
$ws.Range("H2").Value = "KH"
$ws.Range("I2").Value = "TH"
$ws.Range("J2").Value = "KH"
$ws.Range("K2").Value = "TH"

# Merge the header cells.
$ws.Range("A1:A2").Merge()
$ws.Range("B1:C1").Merge()
$ws.Range("D1:E1").Merge()
$ws.Range("F1:G1").Merge()
$ws.Range("H1:I1").Merge()
$ws.Range("J1:K1").Merge()

# Styling to match the new look: Times New Roman font on a darker theme fill,
# centered both horizontally and vertically.
$headerRange = $ws.Range("A1:K2")
$headerRange.Font.Name = "Times New Roman"
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.Interior.ThemeColor = 8
$headerRange.Interior.TintAndShade = -0.4

$ws.Range("I3").Select()

# Make the new sheet the active tab.
$ws.Activate()

$wb.Save()
